# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (column G) values for each game row with freshly
# computed strikeout-derived values (s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 3
    9  = 2
    10 = 2
    11 = 1
    12 = 4
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 5
    18 = 0
    19 = 1
    20 = 1
    21 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
